# 15HP03_index.xlsx — interactive visualisations for leaf profiles and dynamics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window: shrink the sheet-tab area (tabRatio 993 -> 500) ---
$win = $excel.ActiveWindow
$win.TabRatio = 500

# --- Header text: "repeat" -> "rep" (column G header, row 1) ---
$ws.Range("G1").Value = "rep"

# --- Current selection moves from A1:H1 to G2 ---
$ws.Range("G2").Select()

# --- Column widths widen (~12% wider, e.g. after a font-metrics change) ---
$ws.Columns.Item(1).ColumnWidth = 12.833333333333332
$ws.Columns.Item(2).ColumnWidth = 13.5
$ws.Columns.Item(3).ColumnWidth = 15.333333333333332
$ws.Columns.Item(4).ColumnWidth = 9.5
$ws.Columns.Item(5).ColumnWidth = 5.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.0
$ws.Columns.Item(7).ColumnWidth = 8.833333333333332
$ws.Columns.Item(8).ColumnWidth = 11.333333333333332
$ws.Columns.Item(9).ColumnWidth = 13.666666666666668

# --- Page setup: keep fit-to-page on a single page, portrait, letter paper ---
$ps = $ws.PageSetup
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1
$ps.PaperSize = 1
$ps.Orientation = 1
$ps.UsePrinterDefaults = $false
